$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns retain literal text formatting (e.g. trailing
# zeros like "495.40" or "1.00") instead of being auto-coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "56.742.73"
$ws.Range("E2").Value = "  +2.28%  "
$ws.Range("D3").Value = "2.500.78"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "495.40"
$ws.Range("E5").Value = "  +2.73%  "
$ws.Range("D6").Value = "153.77"
$ws.Range("E6").Value = "  +9.19%  "
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "0.514"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").Value = "2.516.94"
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("D10").Value = "5.78"
$ws.Range("E10").Value = "  +5.61%  "
$ws.Range("D11").Value = "0.0993"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("E12").Value = "  +2.76%  "
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "2.933.58"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "56.896.65"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("D16").Value = "21.53"
$ws.Range("E16").Value = "  +4.34%  "
$ws.Range("D17").Value = "0.0000137"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "2.515.67"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").Value = "4.56"
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("D20").Value = "10.38"
$ws.Range("E20").Value = "  +3.76%  "
$ws.Range("D21").Value = "324.95"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("D22").Value = "0.996"
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").Value = "5.93"
$ws.Range("E23").Value = "  +4.32%  "
$ws.Range("D24").Value = "59.10"
$ws.Range("E24").Value = "  +2.09%  "
$ws.Range("D25").Value = "0.412"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  -2.06%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").Value = "2.615.10"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").Value = "7.66"
$ws.Range("E29").Value = "  +3.08%  "
$ws.Range("D30").Value = "0.0₃0823"
$ws.Range("E30").Value = "  +3.57%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").Value = "152.43"
$ws.Range("E32").Value = "  +2.37%  "
$ws.Range("D33").Value = "18.43"
$ws.Range("E33").Value = "  +1.54%  "
$ws.Range("E34").Value = "  +3.40%  "
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("E36").Value = "  +4.63%  "
$ws.Range("D37").Value = "3.80"
$ws.Range("E37").Value = "  +2.54%  "
$ws.Range("D38").Value = "0.876"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("E39").Value = "  +6.24%  "
$ws.Range("D40").Value = "34.29"
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "3.53"
$ws.Range("E41").Value = "  +3.31%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "0.0566"
$ws.Range("E42").Value = "  +2.53%  "
$ws.Range("D43").Value = "0.615"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("D44").Value = "0.995"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").Value = "4.96"
$ws.Range("E45").Value = "  +7.34%  "
$ws.Range("D46").Value = "267.75"
$ws.Range("E46").Value = "  +6.15%  "
$ws.Range("D47").Value = "0.0931"
$ws.Range("E47").Value = "  +3.07%  "
$ws.Range("D48").Value = "0.0231"
$ws.Range("E48").Value = "  +3.29%  "
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("D50").Value = "17.97"
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("D51").Value = "1.912.64"
$ws.Range("E51").Value = "  -3.12%  "
